$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the three "ECs" sending-cluster rows (old rows 2-4). This shifts the
# old "FAPs" sending-cluster rows (5-7) up into rows 2-4.
$ws.Rows.Item(2).Resize(3).Delete() | Out-Null

# Refresh the (now-shifted) FAPs rows with the new TPM-derived numbers.
$row1 = @("FAPs","Il27","Il27ra","ECs",3,1,0.9385806666666667,2.815742,1,1,3,1,0.8792186666666666,2.637656,0.5436518920955525,0.5436518920955525,0.8252176423057778,7.426958780752,0.5436518920955525,0.5436518920955525)
$row2 = @("FAPs","Il27","Il27ra","FAPs",3,1,0.9385806666666667,2.815742,1,1,3,1,0.6612263333333334,1.983679,0.4088595486523692,0.4088595486523693,0.6206142527575557,5.585528274818001,0.4088595486523692,0.4088595486523693)
$row3 = @("FAPs","Il27","Il27ra","MuSCs",3,1,0.9385806666666667,2.815742,1,1,2,0.6666666666666666,0.07680066666666667,0.230402,0.04748855925207817,0.04748855925207817,0.07208362092044446,0.6487525882840001,0.04748855925207817,0.04748855925207817)

$data = @($row1, $row2, $row3)

$rows = $data.Count
$cols = $row1.Count
$arr = New-Object 'object[,]' $rows,$cols
for ($r = 0; $r -lt $rows; $r++) {
    for ($c = 0; $c -lt $cols; $c++) {
        $arr[$r,$c] = $data[$r][$c]
    }
}

$ws.Range("A2:T4").Value = $arr
